$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.13452935218811
$ws.Range("B1").Value = 2.288748502731323
$ws.Range("C1").Value = 11.10049343109131
$ws.Range("D1").Value = 2.117898464202881
$ws.Range("E1").Value = 1.276530027389526
